$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 175.0122446297537
$ws.Range("H2").Value = 1.570385316679309
$ws.Range("I2").Value = 180.8354183511697

$ws.Range("G3").Value = 222.093703762316
$ws.Range("H3").Value = 1.313863758585013
$ws.Range("I3").Value = 757.4603053344842

$ws.Range("G4").Value = 251.2343935383294
$ws.Range("H4").Value = 1.064319297938854
$ws.Range("I4").Value = 1725.381238241695

$ws.Range("G5").Value = 249.423522340056
$ws.Range("H5").Value = 0.8292023139423903
$ws.Range("I5").Value = 2804.581349573278

$ws.Range("G6").Value = 259.1901903256966
$ws.Range("H6").Value = 0.6163706244443756
$ws.Range("I6").Value = 3826.863154164561

$ws.Range("G7").Value = 268.8772083745512
$ws.Range("H7").Value = 0.4310773999481527
$ws.Range("I7").Value = 4774.04539166473

$ws.Range("G8").Value = 275.5460361003784
$ws.Range("H8").Value = 0.2767909289983425
$ws.Range("I8").Value = 5651.632708628164

$ws.Range("G9").Value = 285.7700884429001
$ws.Range("H9").Value = 0.1559267540560175
$ws.Range("I9").Value = 6422.261358150013

$ws.Range("G10").Value = 310.5902366509021
$ws.Range("H10").Value = 0.07316208359586752
$ws.Range("I10").Value = 6939.687666984149

$ws.Range("G11").Value = 316.6519150685446
$ws.Range("H11").Value = 0.01760028951844495
$ws.Range("I11").Value = 7574.011062907172

$ws.Range("G12").Value = -8.470200626034966
$ws.Range("I12").Value = 8042.2015705964

$ws.Range("G13").Value = 316.4765292217521
$ws.Range("H13").Value = 0.01761218883289479
$ws.Range("I13").Value = 7578.617293307808

$ws.Range("G14").Value = 310.4525370000481
$ws.Range("H14").Value = 0.07322036407967489
$ws.Range("I14").Value = 6945.375081411192

$ws.Range("G15").Value = 285.8340165441664
$ws.Range("H15").Value = 0.1560769887523443
$ws.Range("I15").Value = 6427.390927548613

$ws.Range("G16").Value = 275.5382703219476
$ws.Range("H16").Value = 0.2770803374920145
$ws.Range("I16").Value = 5655.971741422171

$ws.Range("G17").Value = 268.9338646196583
$ws.Range("H17").Value = 0.4315482463607944
$ws.Range("I17").Value = 4777.752341287575

$ws.Range("G18").Value = 259.2332626855415
$ws.Range("H18").Value = 0.6170623473856787
$ws.Range("I18").Value = 3829.863976280421

$ws.Range("G19").Value = 249.4921662742798
$ws.Range("H19").Value = 0.8301501560123178
$ws.Range("I19").Value = 2806.812004537876

$ws.Range("G20").Value = 251.2638234148226
$ws.Range("H20").Value = 1.065554606772127
$ws.Range("I20").Value = 1726.781617422622

$ws.Range("G21").Value = 222.1462024803066
$ws.Range("H21").Value = 1.315403179278054
$ws.Range("I21").Value = 758.0860525485501

$ws.Range("G22").Value = 175.0694916177185
$ws.Range("H22").Value = 1.572225661875739
$ws.Range("I22").Value = 180.9722551925083
